# Turn File Contents.xlsx - Sheet1 ("Proposed" / Universe section) updates.
# Renames a couple of header cells (NewTurn -> Intel, RaceTurn -> Orders),
# fills in the previously-missing "current value" cells for the
# AllDesigns/AllFleets/AllMinefields rows, tidies the AllStars caption, and
# clarifies the RaceData row in the "Current" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 ("Universe" / Proposed columns): NewTurn -> Intel, RaceTurn -> Orders
$ws.Range("J2").Value = "Intel"
$ws.Range("L2").Value = "Orders"

# Rows 4-6: add the matching "Proposed" column I entries that mirror column H
$ws.Range("I4").Value = "AllDesigns"
$ws.Range("I5").Value = "AllFleets"
$ws.Range("I6").Value = "AllMinefields"

# Row 9: drop the "(position only)" qualifier
$ws.Range("I9").Value = "AllStars"

# Row 39: clarify that this RaceData field is not a RaceData object
$ws.Range("D39").Value = "RaceData (not a RaceData object)"

# Update the saved view state: scrolled down with D40 as the active cell
$ws.Activate()
$ws.Range("D40").Select()
